# Auto-generated edit script implementing the localization-status.xlsx diff
# (swap b3da3f9c / f9023a3c rows 9 & 10, update status for rows 8 & 9 to 'In Translation')
$wb = $excel.ActiveWorkbook

# ---- Cell value updates ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A9").Value = "f9023a3c-875e-4a97-a4b2-fd06bda66208.md"
$ws1.Range("A10").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md"
$ws1.Range("B8").Value = "In Translation"
$ws1.Range("C8").Value = "In Translation"
$ws1.Range("B9").Value = "In Translation"
$ws1.Range("C9").Value = "In Translation"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A9").Value = "f9023a3c-875e-4a97-a4b2-fd06bda66208.md"
$ws2.Range("C9").Value = "f9023a3c-875e-4a97-a4b2-fd06bda66208.bf365ce712290f5860f70f5d6e3e6ea50c431ee5.zh-cn.xlf"
$ws2.Range("A10").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md"
$ws2.Range("C10").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.zh-cn.xlf"
$ws2.Range("B8").Value = "In Translation"
$ws2.Range("B9").Value = "In Translation"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A9").Value = "f9023a3c-875e-4a97-a4b2-fd06bda66208.md"
$ws3.Range("C9").Value = "f9023a3c-875e-4a97-a4b2-fd06bda66208.bf365ce712290f5860f70f5d6e3e6ea50c431ee5.de-de.xlf"
$ws3.Range("A10").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md"
$ws3.Range("C10").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.de-de.xlf"
$ws3.Range("B8").Value = "In Translation"
$ws3.Range("B9").Value = "In Translation"

# ---- Rebuild hyperlinks (delete then re-add) to sync display text with new values ----
# Sheet: Overview
$ws = $wb.Worksheets.Item(1)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c03c744f431f98459839f925c2e458674a5a7903/e2e/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/eacee6dfb9a6c6f6ad9fe9b4ec9b2887f8b143fb/e2e/4019a77a-fe4f-4df9-8651-217f036e3a2d.md", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/eacee6dfb9a6c6f6ad9fe9b4ec9b2887f8b143fb/e2e/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0f52258048af0f05d91054d991dfabdff1957459/e2e/f36d11f1-0539-48c7-b681-e67590560f73.md", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9163a35061629019473ae71cafe30560eacdddd6/e2e/a819cee7-c250-408c-a4ef-b89806cbb22e.md", "", "", "a819cee7-c250-408c-a4ef-b89806cbb22e.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/357aaf8e4123f01797e23432c687422d6aca724a/e2e/11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md", "", "", "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/e2e/2a8b9950-68ff-4a20-8994-7f443c798e66.md", "", "", "2a8b9950-68ff-4a20-8994-7f443c798e66.md")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/3d12f764c5135e228283c5e95c69a0408e693dd5/e2e/b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md", "", "", "f9023a3c-875e-4a97-a4b2-fd06bda66208.md")
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/e2e/f9023a3c-875e-4a97-a4b2-fd06bda66208.md", "", "", "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md")
$ws.Hyperlinks.Add($ws.Range("A11"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/.localization-config", "", "", ".localization-config")

# Sheet: zh-cn
$ws = $wb.Worksheets.Item(2)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c03c744f431f98459839f925c2e458674a5a7903/e2e/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8d8c301a824d19e0b92448d444fab7d15c817db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.zh-cn.xlf", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c606ea29f71e413eba2235e407ba7cab52a13c83/e2e/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c1d8df3d35487e4578f2a3e917a8b94587a05a12/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/mt/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.zh-cn.xlf", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/eacee6dfb9a6c6f6ad9fe9b4ec9b2887f8b143fb/e2e/4019a77a-fe4f-4df9-8651-217f036e3a2d.md", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d71078bcf1e3a344eede74a12c07e66bbc47d776/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.zh-cn.xlf", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/14a983e830d174ced6901f9e689d43c30aacddad/e2e/4019a77a-fe4f-4df9-8651-217f036e3a2d.md", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a1e0459d2e6873e51137b53ce7d79648b3bb758e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.zh-cn.xlf", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/eacee6dfb9a6c6f6ad9fe9b4ec9b2887f8b143fb/e2e/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d71078bcf1e3a344eede74a12c07e66bbc47d776/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.zh-cn.xlf", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/14a983e830d174ced6901f9e689d43c30aacddad/e2e/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a1e0459d2e6873e51137b53ce7d79648b3bb758e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.zh-cn.xlf", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0f52258048af0f05d91054d991dfabdff1957459/e2e/f36d11f1-0539-48c7-b681-e67590560f73.md", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8d8c301a824d19e0b92448d444fab7d15c817db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.zh-cn.xlf", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c606ea29f71e413eba2235e407ba7cab52a13c83/e2e/f36d11f1-0539-48c7-b681-e67590560f73.md", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.md")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c1d8df3d35487e4578f2a3e917a8b94587a05a12/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/mt/f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.zh-cn.xlf", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9163a35061629019473ae71cafe30560eacdddd6/e2e/a819cee7-c250-408c-a4ef-b89806cbb22e.md", "", "", "a819cee7-c250-408c-a4ef-b89806cbb22e.md")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8d8c301a824d19e0b92448d444fab7d15c817db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/a819cee7-c250-408c-a4ef-b89806cbb22e.1bc42e1f00f384da2bca5a9696aca51e98221387.zh-cn.xlf", "", "", "a819cee7-c250-408c-a4ef-b89806cbb22e.1bc42e1f00f384da2bca5a9696aca51e98221387.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/357aaf8e4123f01797e23432c687422d6aca724a/e2e/11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md", "", "", "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32bd37ab7b9d13c351ab4a9c549def7bcef73d2b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/11dcff0e-746e-4bde-b7da-d7fb64b12d1d.34e508a11b7473be85ca32bf46b6b9b53099aebc.zh-cn.xlf", "", "", "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.34e508a11b7473be85ca32bf46b6b9b53099aebc.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/e2e/2a8b9950-68ff-4a20-8994-7f443c798e66.md", "", "", "2a8b9950-68ff-4a20-8994-7f443c798e66.md")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8d8c301a824d19e0b92448d444fab7d15c817db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/2a8b9950-68ff-4a20-8994-7f443c798e66.63f90d6d79ccf30d7cfff0549b9e72aa98fadf0c.zh-cn.xlf", "", "", "2a8b9950-68ff-4a20-8994-7f443c798e66.63f90d6d79ccf30d7cfff0549b9e72aa98fadf0c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/3d12f764c5135e228283c5e95c69a0408e693dd5/e2e/b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md", "", "", "f9023a3c-875e-4a97-a4b2-fd06bda66208.md")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8d8c301a824d19e0b92448d444fab7d15c817db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.zh-cn.xlf", "", "", "f9023a3c-875e-4a97-a4b2-fd06bda66208.bf365ce712290f5860f70f5d6e3e6ea50c431ee5.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/e2e/f9023a3c-875e-4a97-a4b2-fd06bda66208.md", "", "", "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8d8c301a824d19e0b92448d444fab7d15c817db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/f9023a3c-875e-4a97-a4b2-fd06bda66208.bf365ce712290f5860f70f5d6e3e6ea50c431ee5.zh-cn.xlf", "", "", "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A11"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/.localization-config", "", "", ".localization-config")

# Sheet: de-de
$ws = $wb.Worksheets.Item(3)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c03c744f431f98459839f925c2e458674a5a7903/e2e/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6810aa791bb0698057b05c3f72bd0bd0e0214a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.de-de.xlf", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/71c8793074d7476f2195a7cf671f44ce899cbb56/e2e/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e92fea2ef7b9adadaae729842a5ff36146f27b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/mt/1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.de-de.xlf", "", "", "1fbd4690-653c-4303-8e5c-48eb9ef0c6a0.50e1c57fa4ab52f6e7eaaf1a041975c86d5081e5.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/eacee6dfb9a6c6f6ad9fe9b4ec9b2887f8b143fb/e2e/4019a77a-fe4f-4df9-8651-217f036e3a2d.md", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18f39a2c7beb902f3e7a5e994886438a39fe27e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.de-de.xlf", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e8acad2c66a27c697b1380290cc0742e96651556/e2e/4019a77a-fe4f-4df9-8651-217f036e3a2d.md", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/465ad5b1b3d1108d914b3d464e03ca7836a1f591/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.de-de.xlf", "", "", "4019a77a-fe4f-4df9-8651-217f036e3a2d.fb913f5256b692911142a178bc563f420f248f63.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/eacee6dfb9a6c6f6ad9fe9b4ec9b2887f8b143fb/e2e/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18f39a2c7beb902f3e7a5e994886438a39fe27e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.de-de.xlf", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e8acad2c66a27c697b1380290cc0742e96651556/e2e/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/465ad5b1b3d1108d914b3d464e03ca7836a1f591/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.de-de.xlf", "", "", "5a3c6aca-ee1f-4cd3-837a-ca0b7fa02a6c.0bc733abfc8e5aa9e62a58eecaf5cc6c2db5b32f.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0f52258048af0f05d91054d991dfabdff1957459/e2e/f36d11f1-0539-48c7-b681-e67590560f73.md", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6810aa791bb0698057b05c3f72bd0bd0e0214a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.de-de.xlf", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/71c8793074d7476f2195a7cf671f44ce899cbb56/e2e/f36d11f1-0539-48c7-b681-e67590560f73.md", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.md")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e92fea2ef7b9adadaae729842a5ff36146f27b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/mt/f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.de-de.xlf", "", "", "f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9163a35061629019473ae71cafe30560eacdddd6/e2e/a819cee7-c250-408c-a4ef-b89806cbb22e.md", "", "", "a819cee7-c250-408c-a4ef-b89806cbb22e.md")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6810aa791bb0698057b05c3f72bd0bd0e0214a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/a819cee7-c250-408c-a4ef-b89806cbb22e.1bc42e1f00f384da2bca5a9696aca51e98221387.de-de.xlf", "", "", "a819cee7-c250-408c-a4ef-b89806cbb22e.1bc42e1f00f384da2bca5a9696aca51e98221387.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/357aaf8e4123f01797e23432c687422d6aca724a/e2e/11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md", "", "", "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f61274ba961a1bd40d46147f1b461c395890ec3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/11dcff0e-746e-4bde-b7da-d7fb64b12d1d.34e508a11b7473be85ca32bf46b6b9b53099aebc.de-de.xlf", "", "", "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.34e508a11b7473be85ca32bf46b6b9b53099aebc.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/e2e/2a8b9950-68ff-4a20-8994-7f443c798e66.md", "", "", "2a8b9950-68ff-4a20-8994-7f443c798e66.md")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6810aa791bb0698057b05c3f72bd0bd0e0214a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/2a8b9950-68ff-4a20-8994-7f443c798e66.63f90d6d79ccf30d7cfff0549b9e72aa98fadf0c.de-de.xlf", "", "", "2a8b9950-68ff-4a20-8994-7f443c798e66.63f90d6d79ccf30d7cfff0549b9e72aa98fadf0c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/3d12f764c5135e228283c5e95c69a0408e693dd5/e2e/b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md", "", "", "f9023a3c-875e-4a97-a4b2-fd06bda66208.md")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6810aa791bb0698057b05c3f72bd0bd0e0214a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.de-de.xlf", "", "", "f9023a3c-875e-4a97-a4b2-fd06bda66208.bf365ce712290f5860f70f5d6e3e6ea50c431ee5.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/e2e/f9023a3c-875e-4a97-a4b2-fd06bda66208.md", "", "", "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6810aa791bb0698057b05c3f72bd0bd0e0214a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/f9023a3c-875e-4a97-a4b2-fd06bda66208.bf365ce712290f5860f70f5d6e3e6ea50c431ee5.de-de.xlf", "", "", "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A11"), "https://github.com/OpenLocalizationTest/oltest/blob/e29b048651d99cae4abcffbcc58e0e0a6fe62c37/.localization-config", "", "", ".localization-config")

Write-Output "done"